# Rewrite the <p:cNvPr descr="..."> (Shape.AlternativeText) on the gastrectomy /
# esophagectomy diagram <p:pic> shapes from the old CloudFront URLs to local
# "images/<file>.png" paths.
#
# Quirk #1 this script works around: several slides contain TWO <p:pic> shapes
# that both carry the (invalid, but present-in-source) cNvPr id="0". The
# PowerPoint object-model host here resolves Shapes.Item(N).AlternativeText
# (read AND write) by that duplicated id, so it always lands on the first
# matching <p:pic> in document order no matter which Shapes.Item(N) you use.
# To reach the second picture we temporarily swap the pair's document order
# with Shape.ZOrder(msoBringForward) (a pure in-place reorder - it does not
# clone the shape or touch its r:embed relationship), edit the now-reachable
# "first" picture, then swap back.
#
# Quirk #2 this script works around: calling a function with *named*
# parameters (`Foo -bar $x`) corrupts COM object arguments in this host
# (Shapes.Count reads back as 0 inside the function). Positional calls
# (`Foo $x`) work fine, so every function below is called positionally.

$p = $ppt.ActivePresentation

$msoBringForward = 2

function Get-FirstPicIndex($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        if ($slide.Shapes.Item($i).Type -eq 13) {
            return $i
        }
    }
    return -1
}

function Set-FirstPicDescr($slide, $newDescr) {
    $idx = Get-FirstPicIndex $slide
    $slide.Shapes.Item($idx).AlternativeText = $newDescr
}

function Set-OnePictureSlide($slideIndex, $newDescr) {
    $slide = $p.Slides.Item($slideIndex)
    Set-FirstPicDescr $slide $newDescr
}

function Set-TwoPictureSlide($slideIndex, $descr1, $descr2) {
    $slide = $p.Slides.Item($slideIndex)
    $firstIdx = Get-FirstPicIndex $slide

    # Swap the two pictures' document order so the picture that was second
    # becomes reachable as "first".
    $slide.Shapes.Item($firstIdx).ZOrder($msoBringForward)
    Set-FirstPicDescr $slide $descr2

    # Swap back to restore the original order, then fix the first picture.
    $slide.Shapes.Item($firstIdx).ZOrder($msoBringForward)
    Set-FirstPicDescr $slide $descr1
}

# Slides with a single diagram picture.
Set-OnePictureSlide 8 "images/gast_distal_tumor.png"
Set-OnePictureSlide 9 "images/gast_partial.png"
Set-OnePictureSlide 11 "images/gast_distal_gastrectomy.png"
Set-OnePictureSlide 13 "images/gast_body.png"
Set-OnePictureSlide 14 "images/gast_subtotal.png"
Set-OnePictureSlide 16 "images/gast_proximal_tumor.png"
Set-OnePictureSlide 17 "images/gast_total.png"
Set-OnePictureSlide 19 "images/Eso_tumor00_resection2_1600.png"
Set-OnePictureSlide 20 "images/ivor_lewis_simple2_900.png"
Set-OnePictureSlide 21 "images/ivor_lewis_simple2.png"
Set-OnePictureSlide 22 "images/gast_dualtract.png"

# Slides with two diagram pictures sharing a duplicated shape id.
Set-TwoPictureSlide 10 "images/gast_distal_tumor.png" "images/gast_partial.png"
Set-TwoPictureSlide 12 "images/gast_distal_tumor.png" "images/gast_distal_gastrectomy.png"
Set-TwoPictureSlide 15 "images/gast_body.png" "images/gast_subtotal.png"
Set-TwoPictureSlide 18 "images/gast_proximal_tumor.png" "images/gast_total.png"
Set-TwoPictureSlide 23 "images/gast_proximal_tumor.png" "images/gast_dualtract.png"
